$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 140.46666
$ws.Range("I2").Value = 147.35715
$ws.Range("J2").Value = 44
$ws.Range("K2").Value = 147.35715
$ws.Range("L2").Value = 44
$ws.Range("M2").Value = -34.35714999999999
$ws.Range("N2").Value = -270

$ws.Range("H64").Value = 333337150
$ws.Range("J64").Value = 1000000000
$ws.Range("L64").Value = 1000000000
$ws.Range("N64").Value = -1000000496

$ws.Range("H67").Value = 333337150
$ws.Range("J67").Value = 1000000000
$ws.Range("L67").Value = 1000000000
$ws.Range("N67").Value = -1000001716

$ws.Range("H70").Value = 5436.5
$ws.Range("I70").Value = 2563.5
$ws.Range("K70").Value = 7690.5
$ws.Range("M70").Value = -7420.5

$ws.Range("H73").Value = 5436.5
$ws.Range("I73").Value = 2563.5
$ws.Range("K73").Value = 7690.5
$ws.Range("M73").Value = -6754.5

$ws.Range("H94").Value = 3097.75
$ws.Range("I94").Value = 3097.75
$ws.Range("K94").Value = 3097.75
$ws.Range("M94").Value = -2646.75

$ws.Range("H96").Value = 496.66666
$ws.Range("I96").Value = 270.2
$ws.Range("K96").Value = 810.5999999999999
$ws.Range("M96").Value = 562.4000000000001

$ws.Range("H127").Value = 933
$ws.Range("I127").Value = 933
$ws.Range("K127").Value = 2799
$ws.Range("M127").Value = 2161

$ws.Range("H132").Value = 4896.5454
$ws.Range("I132").Value = 5123.7744
$ws.Range("K132").Value = 15371.3232
$ws.Range("M132").Value = -12841.3232

$ws.Range("H137").Value = 1852.55
$ws.Range("I137").Value = 1866.3125
$ws.Range("K137").Value = 5598.9375
$ws.Range("M137").Value = -3048.9375

$ws.Range("H138").Value = 4696.5
$ws.Range("I138").Value = 4696.5
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 14089.5
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = -8949.5
$ws.Range("N138").ClearContents()

$ws.Range("H141").Value = 3904.8572
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("N54").ClearContents()

$ws.Range("H61").Value = 3653.3635
$ws.Range("I61").Value = 3062.6667
$ws.Range("J61").Value = 3874.875
$ws.Range("K61").Value = 3062.6667
$ws.Range("L61").Value = 3874.875
$ws.Range("M61").Value = -2850.6667
$ws.Range("N61").Value = -4298.875

$ws.Range("H74").Value = 2472.2856
$ws.Range("I74").Value = 1439.8
$ws.Range("J74").Value = 3045.889
$ws.Range("K74").Value = 1439.8
$ws.Range("L74").Value = 3045.889
$ws.Range("M74").Value = -565.8
$ws.Range("N74").Value = -4793.889

$ws.Range("H77").Value = 2472.2856
$ws.Range("I77").Value = 1439.8
$ws.Range("J77").Value = 3045.889
$ws.Range("K77").Value = 7199
$ws.Range("L77").Value = 15229.445
$ws.Range("M77").Value = -2831
$ws.Range("N77").Value = -23965.445

$ws.Range("H105").Value = 99000
$ws.Range("J105").Value = 99000
$ws.Range("L105").Value = 99000
$ws.Range("N105").Value = -105988

$ws.Range("H136").Value = 3653.3635
$ws.Range("I136").Value = 3062.6667
$ws.Range("J136").Value = 3874.875
$ws.Range("K136").Value = 9188.000100000001
$ws.Range("L136").Value = 11624.625
$ws.Range("M136").Value = -6638.000100000001
$ws.Range("N136").Value = -16724.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3277.75
$ws.Range("I86").Value = 3220.125
$ws.Range("K86").Value = 3220.125
$ws.Range("M86").Value = -2097.125

$ws.Range("H89").Value = 3277.75
$ws.Range("I89").Value = 3220.125
$ws.Range("K89").Value = 16100.625
$ws.Range("M89").Value = -10484.625

$ws.Range("I94").Value = 250012820
$ws.Range("K94").Value = 250012820
$ws.Range("M94").Value = -250012369

$ws.Range("H105").Value = 20003084
$ws.Range("I105").Value = 1003079.9
$ws.Range("J105").Value = 83336430
$ws.Range("K105").Value = 1003079.9
$ws.Range("L105").Value = 83336430
$ws.Range("M105").Value = -1001332.9
$ws.Range("N105").Value = -83339924

$ws.Range("H107").Value = 4274721.5
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()

$ws.Range("H134").Value = 1968.5834
$ws.Range("I134").Value = 1250.2
$ws.Range("J134").Value = 3165.889
$ws.Range("K134").Value = 3750.6
$ws.Range("L134").Value = 9497.667000000001
$ws.Range("M134").Value = -1215.6
$ws.Range("N134").Value = -14567.667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5955270.5
$ws.Range("I31").Value = 1924.4
$ws.Range("K31").Value = 1924.4
$ws.Range("M31").Value = -1629.4

$ws.Range("H34").Value = 5955270.5
$ws.Range("I34").Value = 1924.4
$ws.Range("K34").Value = 1924.4
$ws.Range("M34").Value = -1722.4

$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("M42").ClearContents()

$ws.Range("H107").Value = 689.3889
$ws.Range("I107").Value = 471.33334
$ws.Range("J107").Value = 1779.6666
$ws.Range("K107").Value = 471.33334
$ws.Range("L107").Value = 1779.6666
$ws.Range("M107").Value = 1448.66666
$ws.Range("N107").Value = -5619.6666

$ws.Range("H132").Value = 4000.8276
$ws.Range("I132").Value = 3110.7144
$ws.Range("K132").Value = 9332.143199999999
$ws.Range("M132").Value = -6802.143199999999

$ws.Range("H134").Value = 3789.0833
$ws.Range("I134").Value = 3832.5356
$ws.Range("J134").Value = 3637
$ws.Range("K134").Value = 11497.6068
$ws.Range("L134").Value = 10911
$ws.Range("M134").Value = -8962.606800000001
$ws.Range("N134").Value = -15981

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 23749.5
$ws.Range("J23").Value = 23749.5
$ws.Range("L23").Value = 71248.5
$ws.Range("N23").Value = -71718.5

$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()

$ws.Range("H141").Value = 27856.285
$ws.Range("I141").Value = 22998.8
$ws.Range("K141").Value = 68996.39999999999
$ws.Range("M141").Value = -63816.39999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 2966.75
$ws.Range("I2").Value = 2966.75
$ws.Range("K2").Value = 2966.75
$ws.Range("M2").Value = -2853.75

$ws.Range("H80").Value = 125003060
$ws.Range("I80").Value = 200002600
$ws.Range("K80").Value = 200002600
$ws.Range("M80").Value = -200001602

$ws.Range("H83").Value = 125003060
$ws.Range("I83").Value = 200002600
$ws.Range("K83").Value = 1000013000
$ws.Range("M83").Value = -1000008008

$ws.Range("H105").Value = 80000
$ws.Range("J105").Value = 80000
$ws.Range("L105").Value = 80000
$ws.Range("N105").Value = -86988

$ws.Range("H107").Value = 1926.4166
$ws.Range("J107").Value = 6499
$ws.Range("L107").Value = 6499
$ws.Range("N107").Value = -10339

$ws.Range("H118").Value = 46665.332
$ws.Range("J118").Value = 46665.332
$ws.Range("L118").Value = 46665.332
$ws.Range("N118").Value = -49979.332

$ws.Range("H132").Value = 2068.558
$ws.Range("I132").Value = 1940.4193
$ws.Range("K132").Value = 5821.257900000001
$ws.Range("M132").Value = -3291.257900000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3844.2
$ws.Range("I132").Value = 3657.3845
$ws.Range("J132").Value = 4191.143
$ws.Range("K132").Value = 10972.1535
$ws.Range("L132").Value = 12573.429
$ws.Range("M132").Value = -8442.1535
$ws.Range("N132").Value = -17633.429

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 563.5294
$ws.Range("I107").Value = 554.9231
$ws.Range("K107").Value = 1664.7693
$ws.Range("M107").Value = 255.2307000000001

$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
